$wb = $excel.ActiveWorkbook

# --- Tax sheet (sheet2): add GST and Service Tax rows ---
$wsTax = $wb.Worksheets.Item("Tax")

$wsTax.Range("A2").Value = "'9cdbd0c4-8435-4848-9732-c6545ee41faa"
$wsTax.Range("B2").Value = "'GST"
$wsTax.Range("C2").Value = "'15.5"
$wsTax.Range("D2").Value = "'Active"
$wsTax.Range("E2").Value = "'5/8/2019 8:53:31 PM"

$wsTax.Range("A3").Value = "'9aae78bb-0d58-49cb-b91e-25b5d43d60c1"
$wsTax.Range("B3").Value = "'Service Tax"
$wsTax.Range("C3").Value = "'3.5"
$wsTax.Range("D3").Value = "'Active"
$wsTax.Range("E3").Value = "'5/8/2019 8:54:03 PM"

$wsTax.Range("A1:E1").Select() | Out-Null

# --- New ProductType sheet, placed after Tax ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsType = $wb.Worksheets.Add($null, $lastSheet)
$wsType.Name = "ProductType"

$wsType.Range("A1").Value = "ID"
$wsType.Range("B1").Value = "ProductID"
$wsType.Range("C1").Value = "Key"
$wsType.Range("D1").Value = "Value"
$wsType.Range("E1").Value = "Status"
$wsType.Range("F1").Value = "EntryDate"

$wsType.Range("D4").Select() | Out-Null

$wsType.Activate()
